$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 95.57143000000001
$ws.Range("I11").Value = 95.57143000000001
$ws.Range("K11").Value = 95.57143000000001
$ws.Range("M11").Value = 44.42856999999999
$ws.Range("H12").Value = 400.5
$ws.Range("I12").Value = 400.5
$ws.Range("K12").Value = 400.5
$ws.Range("M12").Value = -230.5
$ws.Range("H19").Value = 1198.5
$ws.Range("I19").Value = 1100
$ws.Range("J19").Value = 1231.3334
$ws.Range("K19").Value = 1100
$ws.Range("L19").Value = 1231.3334
$ws.Range("M19").Value = -925
$ws.Range("N19").Value = -1581.3334
$ws.Range("H54").Value = 14998.75
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 14998.75
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 14998.75
$ws.Range("M54").ClearContents()
$ws.Range("N54").Value = -15970.75
$ws.Range("H76").Value = 1000
$ws.Range("J76").Value = 1000
$ws.Range("L76").Value = 1000
$ws.Range("N76").Value = -1630
$ws.Range("H79").Value = 1000
$ws.Range("J79").Value = 1000
$ws.Range("L79").Value = 1000
$ws.Range("N79").Value = -3184
$ws.Range("H80").Value = 528.6667
$ws.Range("J80").Value = 509.66666
$ws.Range("L80").Value = 1528.99998
$ws.Range("N80").Value = -3524.99998
$ws.Range("H83").Value = 528.6667
$ws.Range("J83").Value = 509.66666
$ws.Range("L83").Value = 4586.99994
$ws.Range("N83").Value = -14570.99994
$ws.Range("H86").Value = 1910.25
$ws.Range("J86").Value = 2119.6
$ws.Range("L86").Value = 2119.6
$ws.Range("N86").Value = -4365.6
$ws.Range("H89").Value = 1910.25
$ws.Range("J89").Value = 2119.6
$ws.Range("L89").Value = 10598
$ws.Range("N89").Value = -21830
$ws.Range("H106").Value = 0
$ws.Range("I106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("M106").ClearContents()
$ws.Range("H113").Value = 4212.143
$ws.Range("J113").Value = 4164.1665
$ws.Range("L113").Value = 4164.1665
$ws.Range("N113").Value = -10672.1665
$ws.Range("H116").Value = 5449.75
$ws.Range("H135").Value = 6375.5
$ws.Range("I135").Value = 4865.5
$ws.Range("K135").Value = 43789.5
$ws.Range("M135").Value = -41254.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 22
$ws.Range("I5").Value = 22
$ws.Range("K5").Value = 22
$ws.Range("M5").Value = 90
$ws.Range("H61").Value = 1679.875
$ws.Range("I61").Value = 1746.4286
$ws.Range("J61").Value = 1214
$ws.Range("K61").Value = 1746.4286
$ws.Range("L61").Value = 1214
$ws.Range("M61").Value = -1534.4286
$ws.Range("N61").Value = -1638
$ws.Range("H102").Value = 2278.5557
$ws.Range("I102").Value = 1938.375
$ws.Range("J102").Value = 5000
$ws.Range("K102").Value = 1938.375
$ws.Range("L102").Value = 5000
$ws.Range("M102").Value = -316.375
$ws.Range("N102").Value = -8244
$ws.Range("H122").Value = 1439.2
$ws.Range("I122").Value = 1432.6666
$ws.Range("K122").Value = 4297.9998
$ws.Range("M122").Value = -1847.9998
$ws.Range("H136").Value = 1679.875
$ws.Range("I136").Value = 1746.4286
$ws.Range("J136").Value = 1214
$ws.Range("K136").Value = 5239.2858
$ws.Range("L136").Value = 3642
$ws.Range("M136").Value = -2689.2858
$ws.Range("N136").Value = -8742

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 22
$ws.Range("I4").Value = 22
$ws.Range("K4").Value = 22
$ws.Range("M4").Value = 93
$ws.Range("H20").Value = 5331.3335
$ws.Range("I20").Value = 5331.3335
$ws.Range("K20").Value = 5331.3335
$ws.Range("M20").Value = -5084.3335
$ws.Range("H94").Value = 530.1429000000001
$ws.Range("I94").Value = 582.4
$ws.Range("J94").Value = 399.5
$ws.Range("K94").Value = 582.4
$ws.Range("L94").Value = 399.5
$ws.Range("M94").Value = -131.4
$ws.Range("N94").Value = -1301.5
$ws.Range("H105").Value = 1942.1428
$ws.Range("I105").Value = 1865.8334
$ws.Range("K105").Value = 1865.8334
$ws.Range("M105").Value = -118.8334
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 94.125
$ws.Range("I7").Value = 104.85714
$ws.Range("J7").Value = 19
$ws.Range("K7").Value = 104.85714
$ws.Range("L7").Value = 19
$ws.Range("M7").Value = 8.142859999999999
$ws.Range("N7").Value = -245
$ws.Range("H11").Value = 1641.1111
$ws.Range("I11").Value = 1800
$ws.Range("J11").Value = 1595.7142
$ws.Range("K11").Value = 1800
$ws.Range("L11").Value = 1595.7142
$ws.Range("M11").Value = -1660
$ws.Range("N11").Value = -1875.7142
$ws.Range("H31").Value = 2759.4546
$ws.Range("J31").Value = 3996.25
$ws.Range("L31").Value = 3996.25
$ws.Range("N31").Value = -4586.25
$ws.Range("H34").Value = 2759.4546
$ws.Range("J34").Value = 3996.25
$ws.Range("L34").Value = 3996.25
$ws.Range("N34").Value = -4400.25
$ws.Range("H58").Value = 2326
$ws.Range("J58").Value = 3333
$ws.Range("L58").Value = 3333
$ws.Range("N58").Value = -3739
$ws.Range("H99").Value = 2052.7144
$ws.Range("I99").Value = 1379.5
$ws.Range("J99").Value = 2950.3333
$ws.Range("K99").Value = 1379.5
$ws.Range("L99").Value = 2950.3333
$ws.Range("M99").Value = 118.5
$ws.Range("N99").Value = -5946.3333
$ws.Range("H126").Value = 2052.7144
$ws.Range("I126").Value = 1379.5
$ws.Range("J126").Value = 2950.3333
$ws.Range("K126").Value = 4138.5
$ws.Range("L126").Value = 8850.999899999999
$ws.Range("M126").Value = -1668.5
$ws.Range("N126").Value = -13790.9999
$ws.Range("H132").Value = 5743.75
$ws.Range("J132").Value = 3998.5
$ws.Range("L132").Value = 11995.5
$ws.Range("N132").Value = -17055.5
$ws.Range("H136").Value = 2326
$ws.Range("J136").Value = 3333
$ws.Range("L136").Value = 9999
$ws.Range("N136").Value = -15099
$ws.Range("H141").Value = 38063.715
$ws.Range("J141").Value = 38063.715
$ws.Range("L141").Value = 38063.715
$ws.Range("N141").Value = -48423.715

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 597.1429000000001
$ws.Range("I40").Value = 73.333336
$ws.Range("K40").Value = 293.333344
$ws.Range("M40").Value = -224.333344
$ws.Range("H128").Value = 540938.5
$ws.Range("I128").Value = 540938.5
$ws.Range("K128").Value = 1622815.5
$ws.Range("M128").Value = -1617835.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1436.75
$ws.Range("I80").Value = 1515.6666
$ws.Range("J80").Value = 1200
$ws.Range("K80").Value = 1515.6666
$ws.Range("L80").Value = 1200
$ws.Range("M80").Value = -517.6666
$ws.Range("N80").Value = -3196
$ws.Range("H83").Value = 1436.75
$ws.Range("I83").Value = 1515.6666
$ws.Range("J83").Value = 1200
$ws.Range("K83").Value = 7578.333000000001
$ws.Range("L83").Value = 6000
$ws.Range("M83").Value = -2586.333000000001
$ws.Range("N83").Value = -15984
$ws.Range("H97").Value = 188
$ws.Range("I97").Value = 141.25
$ws.Range("K97").Value = 141.25
$ws.Range("M97").Value = 354.75
$ws.Range("H132").Value = 9245.214
$ws.Range("I132").Value = 9245.214
$ws.Range("K132").Value = 27735.642
$ws.Range("M132").Value = -25205.642

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3965.3333
$ws.Range("I22").Value = 4200
$ws.Range("J22").Value = 3496
$ws.Range("K22").Value = 4200
$ws.Range("L22").Value = 3496
$ws.Range("M22").Value = -3905
$ws.Range("N22").Value = -4086
$ws.Range("H27").Value = 3965.3333
$ws.Range("I27").Value = 4200
$ws.Range("J27").Value = 3496
$ws.Range("K27").Value = 4200
$ws.Range("L27").Value = 3496
$ws.Range("M27").Value = -4093
$ws.Range("N27").Value = -3710
$ws.Range("H55").Value = 1059.1538
$ws.Range("I55").Value = 605.6
$ws.Range("J55").Value = 1342.625
$ws.Range("K55").Value = 605.6
$ws.Range("L55").Value = 1342.625
$ws.Range("M55").Value = -432.6
$ws.Range("N55").Value = -1688.625
$ws.Range("H93").Value = 2227
$ws.Range("I93").Value = 2018
$ws.Range("K93").Value = 2018
$ws.Range("M93").Value = -770
$ws.Range("H132").Value = 3916.5
$ws.Range("I132").Value = 3874.75
$ws.Range("K132").Value = 11624.25
$ws.Range("M132").Value = -9094.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1794.7333
$ws.Range("I107").Value = 1655.4615
$ws.Range("K107").Value = 4966.3845
$ws.Range("M107").Value = -3046.3845
$ws.Range("H113").Value = 548.3333
$ws.Range("I113").Value = 451
$ws.Range("J113").Value = 743
$ws.Range("K113").Value = 1353
$ws.Range("L113").Value = 2229
$ws.Range("M113").Value = 817
$ws.Range("N113").Value = -6569
$ws.Range("H132").Value = 2256.9285
$ws.Range("I132").Value = 2258.0833
$ws.Range("K132").Value = 6774.249899999999
$ws.Range("M132").Value = -4244.249899999999
$ws.Range("H136").Value = 2599
$ws.Range("I136").Value = 2340.3333
$ws.Range("J136").Value = 3375
$ws.Range("K136").Value = 7020.999899999999
$ws.Range("L136").Value = 10125
$ws.Range("M136").Value = -4470.999899999999
$ws.Range("N136").Value = -15225
